$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coste Hora")

# Update the base monthly salary value (B9) from 3000 to 1500.
# All dependent formulas (D9, B10, D10, B11, D11, D17, E17, B27, B28) will
# recalculate automatically.
$ws.Range("B9").Value = 1500

# Update the active selection to match the target state.
$ws.Range("B10").Select()

$wb.Save()
